$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New card list for "Duel Decks Elves vs. Goblins (DD1)"
$values = @(
    'Duel Decks Elves vs. Goblins (DD1)',
    'Akki Coalflinger',
    'Allosaurus Rider',
    'Ambush Commander',
    'Boggart Shenanigans',
    'Clickslither',
    'Elvish Eulogist',
    'Elvish Harbinger',
    'Elvish Promenade',
    'Elvish Warrior',
    'Emberwilde Augur',
    'Flamewave Invoker',
    'Forest',
    'Forest',
    'Forest',
    'Forest',
    'Forgotten Cave',
    'Gempalm Incinerator',
    'Gempalm Strider',
    'Giant Growth',
    'Goblin Burrows',
    'Goblin Cohort',
    'Goblin Matron',
    'Goblin Ringleader',
    'Goblin Sledder',
    'Goblin Warchief',
    'Harmonize',
    'Heedless One',
    'Ib Halfheart, Goblin Tactician',
    'Imperious Perfect',
    'Llanowar Elves',
    'Lys Alana Huntmaster',
    'Mogg Fanatic',
    'Mogg War Marshal',
    'Moonglove Extract',
    'Mountain',
    'Mountain',
    'Mountain',
    'Mountain',
    'Mudbutton Torchrunner',
    'Raging Goblin',
    'Reckless One',
    'Siege-Gang Commander',
    'Skirk Drill Sergeant',
    'Skirk Fire Marshal',
    'Skirk Prospector',
    'Skirk Shaman',
    'Slate of Ancestry',
    'Spitting Earth',
    'Stonewood Invoker',
    'Sylvan Messenger',
    'Tarfire',
    'Tar Pitcher',
    'Timberwatch Elf',
    'Tranquil Thicket',
    'Voice of the Woods',
    'Wellwisher',
    'Wildsize',
    'Wirewood Herald',
    'Wirewood Lodge',
    'Wirewood Symbiote',
    'Wood Elves',
    'Wren''s Run Vanquisher'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# The new sheet has one fewer row (A1:A63) than before (A1:A64); remove the
# now-unused last row.
$ws.Rows.Item(64).Delete()
